# "Forgot to save the spreadsheet before the previous commit"
#
# Adds a new "From Facilities" sheet (a sibling of the existing "From ..."
# distance/time sheets) at the end of the workbook, patterned after the
# last existing sheet ("From Ballard") so it inherits the same layout,
# column widths, cell styles and the shared H:K "minutes" formulas.

$wb = $excel.ActiveWorkbook

# The existing sheets all follow the same template - copy the last one
# ("From Ballard") to the end of the workbook so the new sheet starts out
# with identical formatting/formulas, then overwrite the bits that differ.
$templateSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$templateSheet.Copy($null, $templateSheet)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "From Facilities"

# Header labels above the distance / time tables.
$newSheet.Range("A3").Value = "Distances (feet) from Facilities to…"
$newSheet.Range("G3").Value = "Times (minutes) from Facilities to…"

# Distance data (feet) for each destination - lives in column B for every
# row except the copied template's "self" row, which had its value sitting
# in column C (row 4 / Johnston here); clear that leftover value out.
$newSheet.Range("C4").ClearContents()

$newSheet.Range("B4").Value = 576
$newSheet.Range("B5").Value = 837
$newSheet.Range("B6").Value = 777
$newSheet.Range("B7").Value = 462
$newSheet.Range("B8").Value = 1156
$newSheet.Range("B9").Value = 1784
$newSheet.Range("B10").Value = 1157
$newSheet.Range("B11").Value = 1038
$newSheet.Range("B12").Value = 829
$newSheet.Range("B13").Value = 750
$newSheet.Range("B14").Value = 347
$newSheet.Range("B15").Value = 593
$newSheet.Range("B16").Value = 1409

# Matches the saved selection on the new (now active) sheet.
$newSheet.Range("E25").Select()
